$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Run" flag for the CreateAccount scenario (row 2, column B) from "No" to "Yes"
$ws.Range("B2").Value = "Yes"

# Select the entire 4th row (mirrors clicking the row 4 header in the UI)
$ws.Range("A4:XFD1048576").Select()
